$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 117.42857
$ws.Range("I39").Value = 103.666664
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 310.999992
$ws.Range("L39").Value = 600
$ws.Range("M39").Value = -14.99999200000002
$ws.Range("N39").Value = -1192

# Row 105
$ws.Range("H105").Value = 68717.5
$ws.Range("J105").Value = 67435
$ws.Range("L105").Value = 67435
$ws.Range("N105").Value = -74423

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 129
$ws.Range("H129").Value = 4142.25
$ws.Range("I129").Value = 2248.3333
$ws.Range("J129").Value = 6036.1665
$ws.Range("K129").Value = 6744.999899999999
$ws.Range("L129").Value = 18108.4995
$ws.Range("M129").Value = -1744.999899999999
$ws.Range("N129").Value = -28108.4995

# Row 137
$ws.Range("H137").Value = 2385.2
$ws.Range("I137").Value = 1023.4737
$ws.Range("J137").Value = 4737.273
$ws.Range("K137").Value = 3070.4211
$ws.Range("L137").Value = 14211.819
$ws.Range("M137").Value = -520.4211
$ws.Range("N137").Value = -19311.819

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2440
$ws.Range("I2").Value = 2184
$ws.Range("K2").Value = 2184
$ws.Range("M2").Value = -2071

# Row 32
$ws.Range("H32").Value = 10390.059
$ws.Range("I32").Value = 6603.185
$ws.Range("J32").Value = 24996.572
$ws.Range("K32").Value = 6603.185
$ws.Range("L32").Value = 24996.572
$ws.Range("M32").Value = -6316.185
$ws.Range("N32").Value = -25570.572

# Row 61
$ws.Range("H61").Value = 3096.652
$ws.Range("J61").Value = 4170.6665
$ws.Range("L61").Value = 4170.6665
$ws.Range("N61").Value = -4594.6665

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 116
$ws.Range("H116").Value = 2440
$ws.Range("I116").Value = 2184
$ws.Range("K116").Value = 2184
$ws.Range("M116").Value = 110

# Row 132
$ws.Range("H132").Value = 1599.3334
$ws.Range("J132").Value = 1200
$ws.Range("L132").Value = 3600
$ws.Range("N132").Value = -8660

# Row 136
$ws.Range("H136").Value = 3096.652
$ws.Range("J136").Value = 4170.6665
$ws.Range("L136").Value = 12511.9995
$ws.Range("N136").Value = -17611.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2440
$ws.Range("I3").Value = 2184
$ws.Range("K3").Value = 2184
$ws.Range("M3").Value = -2070

# Row 20
$ws.Range("H20").Value = 6055.2856
$ws.Range("I20").Value = 7679.4
$ws.Range("J20").Value = 1995
$ws.Range("K20").Value = 7679.4
$ws.Range("L20").Value = 1995
$ws.Range("M20").Value = -7432.4
$ws.Range("N20").Value = -2489

# Row 86
$ws.Range("H86").Value = 3903.625
$ws.Range("I86").Value = 3488.8572
$ws.Range("K86").Value = 3488.8572
$ws.Range("M86").Value = -2365.8572

# Row 89
$ws.Range("H89").Value = 3903.625
$ws.Range("I89").Value = 3488.8572
$ws.Range("K89").Value = 17444.286
$ws.Range("M89").Value = -11828.286

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4642.8945
$ws.Range("I31").Value = 1521
$ws.Range("J31").Value = 5757.857
$ws.Range("K31").Value = 1521
$ws.Range("L31").Value = 5757.857
$ws.Range("M31").Value = -1226
$ws.Range("N31").Value = -6347.857

# Row 34
$ws.Range("H34").Value = 4642.8945
$ws.Range("I34").Value = 1521
$ws.Range("J34").Value = 5757.857
$ws.Range("K34").Value = 1521
$ws.Range("L34").Value = 5757.857
$ws.Range("M34").Value = -1319
$ws.Range("N34").Value = -6161.857

# Row 41
$ws.Range("H41").Value = 834
$ws.Range("I41").Value = 834
$ws.Range("K41").Value = 834
$ws.Range("M41").Value = -406

# Row 107
$ws.Range("H107").Value = 878.43475
$ws.Range("I107").Value = 491.08334
$ws.Range("J107").Value = 1301
$ws.Range("K107").Value = 491.08334
$ws.Range("L107").Value = 1301
$ws.Range("M107").Value = 1428.91666
$ws.Range("N107").Value = -5141

# Row 118
$ws.Range("H118").Value = 60000
$ws.Range("J118").Value = 60000
$ws.Range("L118").Value = 60000
$ws.Range("N118").Value = -63314

# Row 134
$ws.Range("H134").Value = 2792.25
$ws.Range("I134").Value = 2188.875
$ws.Range("K134").Value = 6566.625
$ws.Range("M134").Value = -4031.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 16666757
$ws.Range("I7").Value = 16666757
$ws.Range("K7").Value = 50000271
$ws.Range("M7").Value = -50000159

# Row 12
$ws.Range("H12").Value = 748.7222
$ws.Range("J12").Value = 752.38464
$ws.Range("L12").Value = 2257.15392
$ws.Range("N12").Value = -2603.15392

# Row 117
$ws.Range("H117").Value = 2289.75
$ws.Range("J117").Value = 3213.6
$ws.Range("L117").Value = 9640.799999999999
$ws.Range("N117").Value = -16524.8

# Row 121
$ws.Range("H121").Value = 770
$ws.Range("J121").Value = 890
$ws.Range("L121").Value = 2670
$ws.Range("N121").Value = -5290

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# Row 137
$ws.Range("H137").Value = 3946.7778
$ws.Range("I137").Value = 2824.182
$ws.Range("J137").Value = 5710.857
$ws.Range("K137").Value = 8472.545999999998
$ws.Range("L137").Value = 17132.571
$ws.Range("M137").Value = -3372.545999999998
$ws.Range("N137").Value = -27332.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 97
$ws.Range("H97").Value = 820.3
$ws.Range("I97").Value = 784
$ws.Range("K97").Value = 784
$ws.Range("M97").Value = -288

# Row 123
$ws.Range("H123").Value = 114166.836
$ws.Range("J123").Value = 114166.836
$ws.Range("L123").Value = 114166.836
$ws.Range("N123").Value = -119066.836

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 6258.0835
$ws.Range("I61").Value = 6815.4443
$ws.Range("J61").Value = 4586
$ws.Range("K61").Value = 6815.4443
$ws.Range("L61").Value = 4586
$ws.Range("M61").Value = -6613.4443
$ws.Range("N61").Value = -4990

# Row 113
$ws.Range("H113").Value = 6258.0835
$ws.Range("I113").Value = 6815.4443
$ws.Range("J113").Value = 4586
$ws.Range("K113").Value = 6815.4443
$ws.Range("L113").Value = 4586
$ws.Range("M113").Value = -4645.4443
$ws.Range("N113").Value = -8926

# Row 133
$ws.Range("H133").Value = 39999
$ws.Range("J133").Value = 39999
$ws.Range("L133").Value = 39999
$ws.Range("N133").Value = -45059

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 21333.334
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 21333.334
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 21333.334
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -22373.334

# Row 104
$ws.Range("H104").Value = 22995.5
$ws.Range("J104").Value = 22995.5
$ws.Range("L104").Value = 22995.5
$ws.Range("N104").Value = -29983.5
